# Apply corrected fitness values to column C ("Fitness") per the SA algorithm
# correction described in the commit message "correction in sa algorithm and
# 746 logs". Column C holds the best-fitness-so-far values for generations
# 0..250 (rows 2..252); the corrected run collapses into five plateaus.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2..5   (Generation 0..3)   -> 8459
$ws.Range("C2:C5").Value = 8459

# Row 6..17  (Generation 4..15)  -> 8189
$ws.Range("C6:C17").Value = 8189

# Row 18..54 (Generation 16..52) -> 7734
$ws.Range("C18:C54").Value = 7734

# Row 55..92 (Generation 53..90) -> 7704
$ws.Range("C55:C92").Value = 7704

# Row 93..252 (Generation 91..250) -> 7310
$ws.Range("C93:C252").Value = 7310
